$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Folha1")

# Row 3 (Guilherme) - role changes from Actor to programador
$ws.Range("E3").Value = "programador"

# Row 5 - replace Yolo with Francisco
$ws.Range("E5").Value = "tester"
$ws.Range("A5").Value = "Francisco"
$ws.Range("C5").Value = "f.silveira@campus.fct.unl.pt"
$ws.Range("D5").Value = 5

$ws.Range("C5").Select()

$wb.Save()
